$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 9).Value = 382330
}

$ws.Range("I14").Select()
